# Update cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.861.28"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.638.43"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.97%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.44%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.28"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5022"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2566"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06367"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.54"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07731"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.643.04"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.75%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.244"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.865.32"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5446"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7879"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.99"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.893.47"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.19"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.377"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.860"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.964"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.882"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.64"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1133"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.78%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.749"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.79%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.241"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04958"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.96%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.260"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.184"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.539"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.54%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.366"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.59%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.622"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.8891"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.94%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5623"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.142.00"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01563"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.005"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.658"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.80"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8038"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.777.09"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈117"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4533"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.65"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.41%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05047"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.34%  "

